$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Apply matching format (yellow fill) to F17 before setting its value,
# mirroring the already-filled status cells F13:F16.
$ws.Range("F13").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write new shared-string values in the order that matches the target
# shared-strings table layout (57: FrenoTombola, 58: MaquinaHorizontal,
# 59: MaquinaAbajo, 60: PuertaCerrada, 61: PuertaAbierta).
$ws.Range("F17").Value = "FrenoTombola"
$ws.Range("C7").Value = "MaquinaHorizontal"
$ws.Range("C8").Value = "MaquinaAbajo"
$ws.Range("C10").Value = "PuertaCerrada"
$ws.Range("C9").Value = "PuertaAbierta"

$ws.Range("E24").Select()
